$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 2: URL_Path changes from /{userId} to /{userId}/{notebookId}
$ws.Range("B2").Value = "/{userId}/{notebookId}"

# --- New rows: Note endpoints (rows 6-8)
$ws.Range("A6").Value = "Note"
$ws.Range("B6").Value = "/{userId}/{noteId}"
$ws.Range("C6").Value = "GET"
$ws.Range("D6").Value = "получить заметку"

$ws.Range("A7").Value = "Note"
$ws.Range("B7").Value = "/{userId}/{noteId}"
$ws.Range("C7").Value = "DELETE"
$ws.Range("D7").Value = "удалить заметку"

$ws.Range("A8").Value = "Note"
$ws.Range("B8").Value = "/{userId}"
$ws.Range("C8").Value = "POST"
$ws.Range("D8").Value = "получить заметки"

# --- New rows: Tag endpoints (rows 9-12)
$ws.Range("A9").Value = "Tag"
$ws.Range("B9").Value = "/{usedId}/{tagId}"
$ws.Range("C9").Value = "GET"
$ws.Range("D9").Value = "получить тэг"

$ws.Range("A10").Value = "Tag"
$ws.Range("B10").Value = "/{usedId}/{tagId}"
$ws.Range("C10").Value = "DELETE"
$ws.Range("D10").Value = "удалить тэг"

$ws.Range("A11").Value = "Tag"
$ws.Range("B11").Value = "/{usedId}"
$ws.Range("C11").Value = "POST"
$ws.Range("D11").Value = "найти тэг"

$ws.Range("A12").Value = "Tag"
$ws.Range("B12").Value = "/{userId}/{notebookId}/{noteId}"
$ws.Range("C12").Value = "GET"
$ws.Range("D12").Value = "получить  тэги заметки"

# --- Apply the same cell style (s=1, centered) as the other data rows to the new rows
$ws.Range("A6:D12").HorizontalAlignment = -4108
$ws.Range("A6:D12").VerticalAlignment = -4108

# --- Column width adjustments (B and D widened to fit new longer content, as Excel's
#     own "best fit" auto-sizing would do after the new rows are added)
$ws.Columns.Item(2).ColumnWidth = 29.45
$ws.Columns.Item(4).ColumnWidth = 21.45

# --- Update selection to match the final state (L9) and dimension naturally follows data
$ws.Range("L9").Select()
